# Design Critic pass - "Clean White" slide 5 metric layout redesign.
# Combine value+unit at 64pt, widen/restyle the caption, delete the old
# duplicate caption shape and shift the right-hand metric column down.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $cand = $slide.Shapes.Item($i)
        if ($cand.Name -eq $name) {
            return $cand
        }
    }
    return $null
}

# EMU -> point helper (PowerPoint COM reports Left/Top/Width/Height in
# points, backed internally by a 32-bit float, so a naive /12700 division
# can truncate one EMU short after the float32 round-trip; nudge forward
# until the conversion lands back exactly on the requested EMU value).
function Truncate-Val($v) {
    if ($v -ge 0) {
        return [math]::Floor($v)
    } else {
        return [math]::Ceiling($v)
    }
}

function EMU($emuValue) {
    $ideal = $emuValue / 12700.0
    $step = 0.0000001
    for ($i = 0; $i -lt 2000; $i++) {
        $candidate = $ideal + ($step * $i)
        $f = [float]$candidate
        $back = [double]$f * 12700.0
        $trunc = Truncate-Val($back)
        if ($trunc -eq $emuValue) {
            return $candidate
        }
    }
    return $ideal
}

# ---------------------------------------------------------------------
# 1) "Text 0" ("260.000") - move/resize the box, keep the big number at
#    64pt (was 96pt) and append a second run "€" in the same style so the
#    value + currency sign live together in one combined headline.
# ---------------------------------------------------------------------
$valueShape = Get-ShapeByName $s "Text 0"
$valueShape.Left = EMU(457200)
$valueShape.Top = EMU(1645920)
$valueShape.Width = EMU(4572000)
$valueShape.Height = EMU(1097280)

$valueRange = $valueShape.TextFrame.TextRange
$valueRange.Font.Size = 64

$beforeLen = $valueRange.Length
$valueRange.InsertAfter("€") | Out-Null
$afterLen = $valueRange.Length
$euroRun = $valueRange.Characters($beforeLen + 1, $afterLen - $beforeLen)
$euroRun.Font.Size = 64
$euroRun.Font.Bold = $true
$euroColor = $euroRun.Font.Color
$euroColor.RGB = 0

# ---------------------------------------------------------------------
# 2) "Text 1" (the old standalone "€") - repurpose in place to become the
#    "jährliche Einsparung" caption: widen to the full column, drop it
#    below the headline, shrink to 24pt and switch to the muted gray.
# ---------------------------------------------------------------------
$captionShape = Get-ShapeByName $s "Text 1"
$captionShape.Left = EMU(457200)
$captionShape.Top = EMU(2834640)
$captionShape.Width = EMU(4572000)
$captionShape.Height = EMU(365760)

$captionRange = $captionShape.TextFrame.TextRange
$captionRange.Text = "jährliche Einsparung"
$captionRange.Font.Size = 24
$captionRange.Font.Bold = $true
$captionColor = $captionRange.Font.Color
$captionColor.RGB = 6184546

# ---------------------------------------------------------------------
# 3) "Text 2" (the redundant old "jährliche Einsparung" caption) is now
#    fully absorbed into the shape above - remove it.
# ---------------------------------------------------------------------
$dupShape = Get-ShapeByName $s "Text 2"
if ($dupShape -ne $null) {
    $dupShape.Delete()
}

# ---------------------------------------------------------------------
# 4) Divider + right-hand metric column: shift everything down to line
#    up with the shorter left column, renumber the display names to
#    close the gap left by the deleted shape.
# ---------------------------------------------------------------------
$divider = Get-ShapeByName $s "Shape 3"
$divider.Left = EMU(5029200)
$divider.Top = EMU(1645920)
$divider.Width = EMU(18288)
$divider.Height = EMU(2286000)
$divider.Name = "Shape 2"

$metric1Value = Get-ShapeByName $s "Text 4"
$metric1Value.Left = EMU(5486400)
$metric1Value.Top = EMU(1645920)
$metric1Value.Width = EMU(2743200)
$metric1Value.Height = EMU(731520)
$metric1Value.Name = "Text 3"

$metric1Label = Get-ShapeByName $s "Text 5"
$metric1Label.Left = EMU(5486400)
$metric1Label.Top = EMU(2286000)
$metric1Label.Width = EMU(2743200)
$metric1Label.Height = EMU(274320)
$metric1Label.Name = "Text 4"

$metric2Value = Get-ShapeByName $s "Text 6"
$metric2Value.Left = EMU(5486400)
$metric2Value.Top = EMU(2926080)
$metric2Value.Width = EMU(2743200)
$metric2Value.Height = EMU(731520)
$metric2Value.Name = "Text 5"

$metric2Label = Get-ShapeByName $s "Text 7"
$metric2Label.Left = EMU(5486400)
$metric2Label.Top = EMU(3566160)
$metric2Label.Width = EMU(2743200)
$metric2Label.Height = EMU(274320)
$metric2Label.Name = "Text 6"
